$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ambito 1 prompts completos - update comparacion de puntajes values

$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 3

$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("G3").Value = 2

$ws.Range("B4").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 3

$ws.Range("B5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("G5").Value = 3

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 4

$ws.Range("B7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 2

$ws.Range("B8").Value = 2
$ws.Range("D8").Value = 3
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 3

$ws.Range("B9").Value = 3
$ws.Range("D9").Value = 3

$ws.Range("B10").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("F10").Value = 3

$ws.Range("B11").Value = 3
$ws.Range("D11").Value = 2
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 4

$ws.Range("B12").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("F12").Value = 2

$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 2
$ws.Range("F13").Value = 2

$ws.Range("B14").Value = 3
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 3

$ws.Range("B15").Value = 3
$ws.Range("G15").Value = 1

$ws.Range("B16").Value = 2
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = 2

$ws.Range("B17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 2

$ws.Range("B18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 2

$ws.Range("B19").Value = 3
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 2

$ws.Range("B20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("F20").Value = 1

$ws.Range("B21").Value = 2
$ws.Range("D21").Value = 3
$ws.Range("G21").Value = 2
